$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Value2 = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "54.332.32"
Set-TextValue $ws.Range("E2") "  +0.65%  "
Set-TextValue $ws.Range("D3") "2.269.88"
Set-TextValue $ws.Range("E3") "  +0.86%  "
Set-TextValue $ws.Range("E4") "  +0.05%  "
Set-TextValue $ws.Range("D5") "500.38"
Set-TextValue $ws.Range("E5") "  +1.73%  "
Set-TextValue $ws.Range("D6") "129.01"
Set-TextValue $ws.Range("E6") "  +1.53%  "
Set-TextValue $ws.Range("D7") "0.998"
Set-TextValue $ws.Range("E7") "  +0.07%  "
Set-TextValue $ws.Range("E8") "  +0.04%  "
Set-TextValue $ws.Range("D9") "0.0956"
Set-TextValue $ws.Range("E9") "  +0.48%  "
Set-TextValue $ws.Range("E10") "  +1.05%  "
Set-TextValue $ws.Range("D11") "0.336"
Set-TextValue $ws.Range("E11") "  +3.61%  "
Set-TextValue $ws.Range("D12") "4.88"
Set-TextValue $ws.Range("E12") "  +5.06%  "
Set-TextValue $ws.Range("D13") "23.16"
Set-TextValue $ws.Range("E13") "  +6.84%  "
Set-TextValue $ws.Range("D14") "2.671.65"
Set-TextValue $ws.Range("E14") "  +0.70%  "
Set-TextValue $ws.Range("D15") "54.326.90"
Set-TextValue $ws.Range("E15") "  +0.77%  "
Set-TextValue $ws.Range("E16") "  +1.01%  "
Set-TextValue $ws.Range("D17") "2.277.28"
Set-TextValue $ws.Range("E17") "  +0.28%  "
Set-TextValue $ws.Range("E18") "  +3.05%  "
Set-TextValue $ws.Range("E19") "  +1.85%  "
Set-TextValue $ws.Range("D20") "304.37"
Set-TextValue $ws.Range("E20") "  +1.97%  "
Set-TextValue $ws.Range("D21") "6.29"
Set-TextValue $ws.Range("E21") "  -1.93%  "
Set-TextValue $ws.Range("E22") "  +0.11%  "
Set-TextValue $ws.Range("D23") "60.11"
Set-TextValue $ws.Range("E23") "  -2.80%  "
Set-TextValue $ws.Range("D24") "0.999"
Set-TextValue $ws.Range("E24") "  -1.80%  "
Set-TextValue $ws.Range("E25") "  +1.22%  "
Set-TextValue $ws.Range("E26") "  +4.85%  "
Set-TextValue $ws.Range("D27") "174.90"
Set-TextValue $ws.Range("E27") "  +5.35%  "
Set-TextValue $ws.Range("D28") "0.0₃0705"
Set-TextValue $ws.Range("E28") "  +3.95%  "
$ws.Range("B29").Value2 = "PancakeSwap"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D29") "1.61"
Set-TextValue $ws.Range("E29") "  +0.98%  "
$ws.Range("B30").Value2 = "Aptos"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D30") "6.00"
Set-TextValue $ws.Range("E30") "  +3.09%  "
Set-TextValue $ws.Range("D31") "1.08"
Set-TextValue $ws.Range("E31") "  +2.28%  "
Set-TextValue $ws.Range("E32") "  +0.01%  "
Set-TextValue $ws.Range("D33") "17.83"
Set-TextValue $ws.Range("E33") "  +1.23%  "
Set-TextValue $ws.Range("D34") "0.995"
Set-TextValue $ws.Range("E34") "  -0.10%  "
Set-TextValue $ws.Range("D35") "0.951"
Set-TextValue $ws.Range("E35") "  +6.84%  "
Set-TextValue $ws.Range("E36") "  +1.97%  "
Set-TextValue $ws.Range("D37") "3.73"
Set-TextValue $ws.Range("E37") "  +1.95%  "
Set-TextValue $ws.Range("D38") "0.375"
Set-TextValue $ws.Range("E38") "  +1.11%  "
Set-TextValue $ws.Range("E39") "  +0.96%  "
Set-TextValue $ws.Range("D40") "3.38"
Set-TextValue $ws.Range("E40") "  +1.32%  "
Set-TextValue $ws.Range("E41") "  -1.76%  "
Set-TextValue $ws.Range("D42") "124.69"
Set-TextValue $ws.Range("E42") "  -0.16%  "
Set-TextValue $ws.Range("D43") "0.0491"
Set-TextValue $ws.Range("E43") "  +2.11%  "
Set-TextValue $ws.Range("D44") "0.0899"
Set-TextValue $ws.Range("E44") "  +1.38%  "
Set-TextValue $ws.Range("D45") "245.78"
Set-TextValue $ws.Range("E45") "  +4.37%  "
Set-TextValue $ws.Range("D46") "0.546"
Set-TextValue $ws.Range("E46") "  +1.01%  "
Set-TextValue $ws.Range("D47") "0.374"
Set-TextValue $ws.Range("E47") "  +1.27%  "
Set-TextValue $ws.Range("D48") "0.0206"
Set-TextValue $ws.Range("E48") "  +2.37%  "
Set-TextValue $ws.Range("E49") "  +0.70%  "
Set-TextValue $ws.Range("E50") "  +1.21%  "
Set-TextValue $ws.Range("E51") "  +3.47%  "
